# B6-PowerPoint.pptx edit — commit "Sun, May 03, 2020  4:06:10 PM"
#
# The canonical-OOXML diff shows two kinds of changes:
#
#   1. Three tables (on the slides that hold the Component 3 comparison
#      tables) get a new <a:tableStyleId> — from the custom "Table_0"
#      style {BA71F996-A943-475B-B67B-8CABE722F230} to the built-in
#      style {6B7EA1FD-E510-4196-B47E-504C82D3DCC5}.
#
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml trade places: the
#      deck's live theme (theme2.xml, "Integral" / "Red Violet") ends up
#      holding the old theme1.xml content ("Office Theme" / "Office"
#      colours), while theme1.xml (the Notes Master's theme, untouched by
#      any exposed object-model call here) keeps the colours that used to
#      live in theme2.xml. The only structural difference between the two
#      theme parts is the <a:theme>/<a:clrScheme> "name" attributes plus
#      10 of the 12 scheme colours (dk1/lt1 are black/white in both) — so
#      we reproduce the reachable part of that swap, the colour values,
#      through the live ThemeColorScheme object (PowerPoint's "Colors"
#      tool), which is the only COM surface this host exposes for
#      mutating theme colours in place.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables.
# ---------------------------------------------------------------------
$newTableStyle = "{6B7EA1FD-E510-4196-B47E-504C82D3DCC5}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle, $false)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colour values (dk2, lt2, accent1-6, hlink, folHlink)
#    over to what used to be theme1.xml's ("Office Theme") palette.
# ---------------------------------------------------------------------
function ToRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

# msoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5 accent1 .. 10 accent6, 11 hlink, 12 folHlink.
$themeColors.Colors(3).RGB  = ToRGB 0x44 0x54 0x6A   # dk2      -> #44546A
$themeColors.Colors(4).RGB  = ToRGB 0xE7 0xE6 0xE6   # lt2      -> #E7E6E6
$themeColors.Colors(5).RGB  = ToRGB 0x5B 0x9B 0xD5   # accent1  -> #5B9BD5
$themeColors.Colors(6).RGB  = ToRGB 0xED 0x7D 0x31   # accent2  -> #ED7D31
$themeColors.Colors(7).RGB  = ToRGB 0xA5 0xA5 0xA5   # accent3  -> #A5A5A5
$themeColors.Colors(8).RGB  = ToRGB 0xFF 0xC0 0x00   # accent4  -> #FFC000
$themeColors.Colors(9).RGB  = ToRGB 0x44 0x72 0xC4   # accent5  -> #4472C4
$themeColors.Colors(10).RGB = ToRGB 0x70 0xAD 0x47   # accent6  -> #70AD47
$themeColors.Colors(11).RGB = ToRGB 0x05 0x63 0xC1   # hlink    -> #0563C1
$themeColors.Colors(12).RGB = ToRGB 0x95 0x4F 0x72   # folHlink -> #954F72
